$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: new Time Recording Log entry (11월 11일 / Initial Data Set 자료 추가) ---
$ws.Range("A31").Value = "11월 11일"
$ws.Range("A31").Characters(3, 1).Font.Name = "돋움"
$ws.Range("A31").Characters(3, 1).Font.ColorIndex = -4105
$ws.Range("A31").Characters(4, 3).Font.Name = "Arial"
$ws.Range("A31").Characters(4, 3).Font.ColorIndex = -4105
$ws.Range("A31").Characters(7, 1).Font.Name = "돋움"
$ws.Range("A31").Characters(7, 1).Font.ColorIndex = -4105

$ws.Range("B31").Value = 0.66666666666666663
$ws.Range("C31").Value = 0.69027777777777777
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 34

$ws.Range("F31").Value = "Initial Data Set 자료 추가"
$ws.Range("F31").Characters(18, 2).Font.Name = "돋움"
$ws.Range("F31").Characters(18, 2).Font.ColorIndex = -4105
$ws.Range("F31").Characters(20, 1).Font.Name = "Arial"
$ws.Range("F31").Characters(20, 1).Font.ColorIndex = -4105
$ws.Range("F31").Characters(21, 2).Font.Name = "돋움"
$ws.Range("F31").Characters(21, 2).Font.ColorIndex = -4105

# --- Row 32: new Time Recording Log entry (11월 12일 / Initial Data Set 자료 추가) ---
$ws.Range("A32").Value = "11월 12일"
$ws.Range("A32").Characters(3, 1).Font.Name = "돋움"
$ws.Range("A32").Characters(3, 1).Font.ColorIndex = -4105
$ws.Range("A32").Characters(4, 3).Font.Name = "Arial"
$ws.Range("A32").Characters(4, 3).Font.ColorIndex = -4105
$ws.Range("A32").Characters(7, 1).Font.Name = "돋움"
$ws.Range("A32").Characters(7, 1).Font.ColorIndex = -4105

$ws.Range("B32").Value = 0.68125000000000002
$ws.Range("C32").Value = 0.74930555555555556
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 98

$ws.Range("F32").Value = "Initial Data Set 자료 추가"

# Move the view's active selection to reflect the newly entered rows
$ws.Range("F33").Select()
